$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 329
$ws.Range("I6").Value = 329
$ws.Range("K6").Value = 987
$ws.Range("M6").Value = -875
$ws.Range("H51").Value = 5626
$ws.Range("I51").Value = 3499
$ws.Range("J51").Value = 5980.5
$ws.Range("K51").Value = 3499
$ws.Range("L51").Value = 5980.5
$ws.Range("M51").Value = -3015
$ws.Range("N51").Value = -6948.5
$ws.Range("H88").Value = 1999999
$ws.Range("I88").Value = 1999999
$ws.Range("K88").Value = 1999999
$ws.Range("M88").Value = -1999593
$ws.Range("H91").Value = 1999999
$ws.Range("I91").Value = 1999999
$ws.Range("K91").Value = 1999999
$ws.Range("M91").Value = -1998595
$ws.Range("H98").Value = 2845.5386
$ws.Range("J98").Value = 2854.7144
$ws.Range("L98").Value = 2854.7144
$ws.Range("N98").Value = -5850.7144
$ws.Range("H116").Value = 5915.9165
$ws.Range("I116").Value = 5815.778
$ws.Range("J116").Value = 6216.3335
$ws.Range("K116").Value = 5815.778
$ws.Range("L116").Value = 6216.3335
$ws.Range("M116").Value = -2373.778
$ws.Range("N116").Value = -13100.3335
$ws.Range("H122").Value = 2845.5386
$ws.Range("J122").Value = 2854.7144
$ws.Range("L122").Value = 8564.143199999999
$ws.Range("N122").Value = -13464.1432
$ws.Range("H125").Value = 5473.3335
$ws.Range("I125").Value = 5215.5
$ws.Range("K125").Value = 46939.5
$ws.Range("M125").Value = -44479.5
$ws.Range("H129").Value = 3281.7144
$ws.Range("I129").Value = 3390.2
$ws.Range("J129").Value = 3010.5
$ws.Range("K129").Value = 10170.6
$ws.Range("L129").Value = 9031.5
$ws.Range("M129").Value = -5170.599999999999
$ws.Range("N129").Value = -19031.5
$ws.Range("H132").Value = 187298.19
$ws.Range("I132").Value = 1988.1957
$ws.Range("K132").Value = 5964.5871
$ws.Range("M132").Value = -3434.5871
$ws.Range("H138").Value = 2347.12
$ws.Range("J138").Value = 2081
$ws.Range("L138").Value = 6243
$ws.Range("N138").Value = -16523

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 365
$ws.Range("I4").Value = 365
$ws.Range("K4").Value = 365
$ws.Range("M4").Value = -249
$ws.Range("H6").Value = 420
$ws.Range("I6").Value = 420
$ws.Range("K6").Value = 420
$ws.Range("M6").Value = -247
$ws.Range("H45").Value = 2388.125
$ws.Range("I45").Value = 2039.4615
$ws.Range("K45").Value = 2039.4615
$ws.Range("M45").Value = -1662.4615
$ws.Range("H88").Value = 2501.25
$ws.Range("J88").Value = 2007
$ws.Range("L88").Value = 2007
$ws.Range("N88").Value = -2819
$ws.Range("H91").Value = 2501.25
$ws.Range("J91").Value = 2007
$ws.Range("L91").Value = 2007
$ws.Range("N91").Value = -4815
$ws.Range("H110").Value = 961.2941
$ws.Range("I110").Value = 862.1667
$ws.Range("J110").Value = 1199.2
$ws.Range("K110").Value = 862.1667
$ws.Range("L110").Value = 1199.2
$ws.Range("M110").Value = 1182.8333
$ws.Range("N110").Value = -5289.2

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 49770.5
$ws.Range("J63").Value = 70271
$ws.Range("L63").Value = 70271
$ws.Range("N63").Value = -71643
$ws.Range("H66").Value = 49770.5
$ws.Range("J66").Value = 70271
$ws.Range("L66").Value = 210813
$ws.Range("N66").Value = -217677
$ws.Range("H94").Value = 1965.65
$ws.Range("I94").Value = 1676
$ws.Range("K94").Value = 1676
$ws.Range("M94").Value = -1225
$ws.Range("H99").Value = 1565.2222
$ws.Range("I99").Value = 1565.2222
$ws.Range("K99").Value = 1565.2222
$ws.Range("M99").Value = -67.22219999999993

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 4838.2856
$ws.Range("J31").Value = 4348.5
$ws.Range("L31").Value = 4348.5
$ws.Range("N31").Value = -4938.5
$ws.Range("H34").Value = 4838.2856
$ws.Range("J34").Value = 4348.5
$ws.Range("L34").Value = 4348.5
$ws.Range("N34").Value = -4752.5
$ws.Range("H64").Value = 29994
$ws.Range("J64").Value = 29994
$ws.Range("L64").Value = 29994
$ws.Range("N64").Value = -30490
$ws.Range("H67").Value = 29994
$ws.Range("J67").Value = 29994
$ws.Range("L67").Value = 29994
$ws.Range("N67").Value = -31710
$ws.Range("H99").Value = 40550.4
$ws.Range("I99").Value = 8188.125
$ws.Range("K99").Value = 8188.125
$ws.Range("M99").Value = -6690.125
$ws.Range("H126").Value = 40550.4
$ws.Range("I126").Value = 8188.125
$ws.Range("K126").Value = 24564.375
$ws.Range("M126").Value = -22094.375
$ws.Range("H134").Value = 2595.1333
$ws.Range("I134").Value = 2852.375
$ws.Range("K134").Value = 8557.125
$ws.Range("M134").Value = -6022.125

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 490.18182
$ws.Range("I29").Value = 115
$ws.Range("J29").Value = 527.7
$ws.Range("K29").Value = 345
$ws.Range("L29").Value = 1583.1
$ws.Range("M29").Value = -68
$ws.Range("N29").Value = -2137.1
$ws.Range("H121").Value = 33534.535
$ws.Range("J121").Value = 296.33334
$ws.Range("L121").Value = 889.0000200000001
$ws.Range("N121").Value = -3509.00002
$ws.Range("H129").Value = 4050.3845
$ws.Range("J129").Value = 4889.5
$ws.Range("L129").Value = 14668.5
$ws.Range("N129").Value = -24668.5
$ws.Range("H132").Value = 2875.8572
$ws.Range("I132").Value = 2159.7334
$ws.Range("J132").Value = 4666.1665
$ws.Range("K132").Value = 19437.6006
$ws.Range("L132").Value = 41995.4985
$ws.Range("M132").Value = -16907.6006
$ws.Range("N132").Value = -47055.4985

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 212.88889
$ws.Range("I2").Value = 225
$ws.Range("J2").Value = 116
$ws.Range("K2").Value = 225
$ws.Range("L2").Value = 116
$ws.Range("M2").Value = -112
$ws.Range("N2").Value = -342
$ws.Range("H34").Value = 29994
$ws.Range("J34").Value = 29994
$ws.Range("L34").Value = 29994
$ws.Range("N34").Value = -30530
$ws.Range("H76").Value = 29994
$ws.Range("J76").Value = 29994
$ws.Range("L76").Value = 29994
$ws.Range("N76").Value = -30624
$ws.Range("H79").Value = 29994
$ws.Range("J79").Value = 29994
$ws.Range("L79").Value = 29994
$ws.Range("N79").Value = -32178
$ws.Range("H102").Value = 2241.423
$ws.Range("I102").Value = 1650.2354
$ws.Range("K102").Value = 1650.2354
$ws.Range("M102").Value = -28.23540000000003
$ws.Range("H126").Value = 4404.143
$ws.Range("I126").Value = 2916.5
$ws.Range("J126").Value = 4999.2
$ws.Range("K126").Value = 8749.5
$ws.Range("L126").Value = 14997.6
$ws.Range("M126").Value = -6279.5
$ws.Range("N126").Value = -19937.6
$ws.Range("H132").Value = 3593.077
$ws.Range("I132").Value = 2233
$ws.Range("K132").Value = 6699
$ws.Range("M132").Value = -4169

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4424.25
$ws.Range("J7").Value = 4999
$ws.Range("L7").Value = 4999
$ws.Range("N7").Value = -5223
$ws.Range("H40").Value = 6887.3335
$ws.Range("I40").Value = 6855.7144
$ws.Range("K40").Value = 6855.7144
$ws.Range("M40").Value = -6719.7144
$ws.Range("H74").Value = 39600
$ws.Range("I74").Value = 39300
$ws.Range("J74").Value = 39675
$ws.Range("K74").Value = 39300
$ws.Range("L74").Value = 39675
$ws.Range("M74").Value = -38302
$ws.Range("N74").Value = -41671
$ws.Range("H76").Value = 31246
$ws.Range("J76").Value = 31246
$ws.Range("L76").Value = 31246
$ws.Range("N76").Value = -31922
$ws.Range("H77").Value = 39600
$ws.Range("I77").Value = 39300
$ws.Range("J77").Value = 39675
$ws.Range("K77").Value = 117900
$ws.Range("L77").Value = 119025
$ws.Range("M77").Value = -112908
$ws.Range("N77").Value = -129009
$ws.Range("H79").Value = 31246
$ws.Range("J79").Value = 31246
$ws.Range("L79").Value = 31246
$ws.Range("N79").Value = -33586
$ws.Range("H122").Value = 6426.5713
$ws.Range("I122").Value = 6163
$ws.Range("J122").Value = 6624.25
$ws.Range("K122").Value = 18489
$ws.Range("L122").Value = 19872.75
$ws.Range("M122").Value = -16039
$ws.Range("N122").Value = -24772.75
$ws.Range("H126").Value = 4424.25
$ws.Range("J126").Value = 4999
$ws.Range("L126").Value = 14997
$ws.Range("N126").Value = -19937

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2810
$ws.Range("I81").Value = 2413.3333
$ws.Range("K81").Value = 4826.6666
$ws.Range("M81").Value = -3765.6666
$ws.Range("H84").Value = 2810
$ws.Range("I84").Value = 2413.3333
$ws.Range("K84").Value = 24133.333
$ws.Range("M84").Value = -18829.333
$ws.Range("H136").Value = 1797.2
$ws.Range("I136").Value = 1591.6666
$ws.Range("J136").Value = 1885.2858
$ws.Range("K136").Value = 4774.9998
$ws.Range("L136").Value = 5655.857400000001
$ws.Range("M136").Value = -2224.9998
$ws.Range("N136").Value = -10755.8574
